$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.026.18'
$ws.Range("E2").Value = '  +2.39%  '
$ws.Range("D3").Value = '1.673.92'
$ws.Range("E3").Value = '  +3.44%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '216.11'
$ws.Range("E5").Value = '  +1.55%  '
$ws.Range("E6").Value = '  +2.06%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  +2.73%  '
$ws.Range("D9").Value = '0.0617'
$ws.Range("E9").Value = '  +1.59%  '
$ws.Range("D10").Value = '20.17'
$ws.Range("E10").Value = '  +5.33%  '
$ws.Range("D11").Value = '0.0887'
$ws.Range("E11").Value = '  +4.73%  '
$ws.Range("D12").Value = '1.909.61'
$ws.Range("E12").Value = '  +3.46%  '
$ws.Range("D13").Value = '1.672.13'
$ws.Range("E13").Value = '  +3.38%  '
$ws.Range("E14").Value = '  +1.76%  '
$ws.Range("E15").Value = '  +2.77%  '
$ws.Range("D16").Value = '65.85'
$ws.Range("E16").Value = '  +3.22%  '
$ws.Range("D17").Value = '27.030.40'
$ws.Range("E17").Value = '  +2.36%  '
$ws.Range("D18").Value = '236.10'
$ws.Range("E18").Value = '  +0.02%  '
$ws.Range("E19").Value = '  +1.77%  '
$ws.Range("D20").Value = '7.73'
$ws.Range("E20").Value = '  -0.81%  '
$ws.Range("E21").Value = '  +0.07%  '
$ws.Range("E22").Value = '  +4.03%  '
$ws.Range("E23").Value = '  +2.07%  '
$ws.Range("E24").Value = '  +2.35%  '
$ws.Range("D25").Value = '145.62'
$ws.Range("E25").Value = '  -0.99%  '
$ws.Range("E26").Value = '  +1.46%  '
$ws.Range("E27").Value = '  +0.61%  '
$ws.Range("D28").Value = '15.95'
$ws.Range("E28").Value = '  +2.59%  '
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("E30").Value = '  +0.38%  '
$ws.Range("E31").Value = '  +1.90%  '
$ws.Range("E32").Value = '  +2.04%  '
$ws.Range("D33").Value = '1.473.74'
$ws.Range("E33").Value = '  -2.83%  '
$ws.Range("E35").Value = '  +6.36%  '
$ws.Range("E36").Value = '  -0.49%  '
$ws.Range("E37").Value = '  +1.29%  '
$ws.Range("D38").Value = '0.897'
$ws.Range("E38").Value = '  +7.61%  '
$ws.Range("E39").Value = '  +2.06%  '
$ws.Range("D40").Value = '6.12'
$ws.Range("E40").Value = '  +3.81%  '
$ws.Range("E41").Value = '  +0.06%  '
$ws.Range("D42").Value = '1.01'
$ws.Range("E42").Value = '  +11.43%  '
$ws.Range("E43").Value = '  +3.50%  '
$ws.Range("D44").Value = '66.52'
$ws.Range("E44").Value = '  +7.54%  '
$ws.Range("D45").Value = '1.818.60'
$ws.Range("E45").Value = '  +3.45%  '
$ws.Range("E46").Value = '  +2.38%  '
$ws.Range("D47").Value = '90.16'
$ws.Range("E47").Value = '  +0.04%  '
$ws.Range("E48").Value = '  +2.04%  '
$ws.Range("E49").Value = '  +4.08%  '
$ws.Range("D50").Value = '0.0507'
$ws.Range("E50").Value = '  +1.19%  '
$ws.Range("D51").Value = '7.63'
$ws.Range("E51").Value = '  +2.11%  '
